$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value looks numeric need to be forced to Text
# so Excel does not convert them to a floating point number (the source
# data is plain text, e.g. "311.85", "0.09541", "93.39", ...).
$textCells = @(
    "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D16", "D17",
    "D18", "D19", "D20", "D22", "D25", "D26", "D28", "D29", "D30", "D31",
    "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44",
    "D45", "D46", "D47", "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values row by row
# Row 2
$ws.Range("D2").Value = "28.219.08"
$ws.Range("E2").Value = "  +0.65%  "

# Row 3
$ws.Range("D3").Value = "1.873.47"
$ws.Range("E3").Value = "  +4.03%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "311.85"
$ws.Range("E5").Value = "  +0.68%  "

# Row 6
$ws.Range("E6").Value = "  -0.11%  "

# Row 7
$ws.Range("D7").Value = "0.5010"
$ws.Range("E7").Value = "  -1.22%  "

# Row 8
$ws.Range("D8").Value = "0.3920"
$ws.Range("E8").Value = "  +2.44%  "

# Row 9
$ws.Range("D9").Value = "0.09541"
$ws.Range("E9").Value = "  +23.56%  "

# Row 10
$ws.Range("D10").Value = "1.142"
$ws.Range("E10").Value = "  +4.52%  "

# Row 11
$ws.Range("D11").Value = "41.01"
$ws.Range("E11").Value = "  +0.73%  "

# Row 12
$ws.Range("D12").Value = "6.480"
$ws.Range("E12").Value = "  +2.10%  "

# Row 13
$ws.Range("D13").Value = "21.00"
$ws.Range("E13").Value = "  +3.36%  "

# Row 14
$ws.Range("D14").Value = "1.877.72"
$ws.Range("E14").Value = "  +4.19%  "

# Row 15
$ws.Range("E15").Value = "  -0.01%  "

# Row 16
$ws.Range("D16").Value = "7.407"
$ws.Range("E16").Value = "  +1.91%  "

# Row 17
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "93.39"
$ws.Range("E17").Value = "  +1.49%  "

# Row 18
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.00001122"
$ws.Range("E18").Value = "  +4.78%  "

# Row 19
$ws.Range("D19").Value = "0.06628"
$ws.Range("E19").Value = "  +0.93%  "

# Row 20
$ws.Range("D20").Value = "17.51"
$ws.Range("E20").Value = "  +1.69%  "

# Row 21
$ws.Range("E21").Value = "  -0.06%  "

# Row 22
$ws.Range("D22").Value = "6.150"
$ws.Range("E22").Value = "  +2.79%  "

# Row 23
$ws.Range("D23").Value = "28.280.06"
$ws.Range("E23").Value = "  +0.80%  "

# Row 24
$ws.Range("E24").Value = "  +2.69%  "

# Row 25
$ws.Range("D25").Value = "2.277"
$ws.Range("E25").Value = "  +2.57%  "

# Row 26
$ws.Range("D26").Value = "2.548"
$ws.Range("E26").Value = "  +5.48%  "

# Row 27
$ws.Range("D27").Value = "2.091.24"
$ws.Range("E27").Value = "  +4.00%  "

# Row 28
$ws.Range("D28").Value = "21.19"
$ws.Range("E28").Value = "  +4.88%  "

# Row 29
$ws.Range("D29").Value = "157.42"
$ws.Range("E29").Value = "  -1.07%  "

# Row 30
$ws.Range("D30").Value = "127.60"
$ws.Range("E30").Value = "  +0.49%  "

# Row 31
$ws.Range("D31").Value = "1.069"
$ws.Range("E31").Value = "  +2.39%  "

# Row 32
$ws.Range("E32").Value = "  -3.39%  "

# Row 33
$ws.Range("D33").Value = "5.642"
$ws.Range("E33").Value = "  +1.93%  "

# Row 34
$ws.Range("E34").Value = "  -0.47%  "

# Row 35
$ws.Range("D35").Value = "0.06756"
$ws.Range("E35").Value = "  -2.68%  "

# Row 36
$ws.Range("D36").Value = "9.536"
$ws.Range("E36").Value = "  +4.99%  "

# Row 37
$ws.Range("D37").Value = "0.02398"
$ws.Range("E37").Value = "  +3.14%  "

# Row 38
$ws.Range("D38").Value = "0.2177"
$ws.Range("E38").Value = "  +0.57%  "

# Row 39
$ws.Range("D39").Value = "11.51"
$ws.Range("E39").Value = "  +0.79%  "

# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.6334"
$ws.Range("E40").Value = "  +3.87%  "

# Row 41
$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").Value = "4.983"
$ws.Range("E41").Value = "  -0.19%  "

# Row 42
$ws.Range("D42").Value = "1.175"
$ws.Range("E42").Value = "  +2.24%  "

# Row 44
$ws.Range("D44").Value = "13.49"
$ws.Range("E44").Value = "  +2.15%  "

# Row 45
$ws.Range("D45").Value = "0.6064"
$ws.Range("E45").Value = "  +3.24%  "

# Row 46
$ws.Range("D46").Value = "3.664"
$ws.Range("E46").Value = "  -1.01%  "

# Row 47
$ws.Range("D47").Value = "1.262"
$ws.Range("E47").Value = "  -2.41%  "

# Row 48
$ws.Range("E48").Value = "  -1.25%  "

# Row 49
$ws.Range("D49").Value = "1.989"
$ws.Range("E49").Value = "  +3.16%  "

# Row 50
$ws.Range("D50").Value = "1.195"
$ws.Range("E50").Value = "  +0.93%  "

# Row 51
$ws.Range("D51").Value = "0.06851"
$ws.Range("E51").Value = "  +1.99%  "
